$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.516.47"
$ws.Range("E2").Value = "  +4.25%  "
$ws.Range("D3").Value = "3.464.14"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("D8").Value = "3.455.93"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.646"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("D15").Value = "3.999.02"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "3.454.70"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "67.425.46"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.120"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.47%  "
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.66%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "595.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.149"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.19%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0786"
$ws.Range("E39").Value = "  +5.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.390"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.41%  "
$ws.Range("D42").Value = "3.144.38"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("E43").Value = "  +4.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.41%  "
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("E46").Value = "  +22.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.16%  "
